$d = $word.ActiveDocument

# Change the three section headings from Title style to Heading1
$d.Paragraphs.Item(1).Style = "Heading1"
$d.Paragraphs.Item(5).Style = "Heading1"
$d.Paragraphs.Item(9).Style = "Heading1"

# Replace the topic-model result text with the new (stemmed / TF-IDF) output
$d.Paragraphs.Item(2).Range.Text = "(0, u'0.000*`"oper`" + 0.000*`"necessari`" + 0.000*`"bike`" + 0.000*`"medium`" + 0.000*`"stabl`" + 0.000*`"winter`" + 0.000*`"core`" + 0.000*`"organis`" + 0.000*`"cost`" + 0.000*`"daili`" + 0.000*`"coordin`" + 0.000*`"local`" + 0.000*`"key`" + 0.000*`"mechan`" + 0.000*`"segment`" + 0.000*`"indirectli`" + 0.000*`"coher`" + 0.000*`"complex`" + 0.000*`"activ`" + 0.000*`"abstract`" + 0.000*`"region`" + 0.000*`"fish`" + 0.000*`"depress`" + 0.000*`"sell`" + 0.000*`"worth`" + 0.000*`"danijel`" + 0.000*`"labour`" + 0.000*`"waitress`" + 0.000*`"cab`" + 0.000*`"alvin`"')"
$d.Paragraphs.Item(3).Range.Text = "(1, u'0.000*`"21st`" + 0.000*`"hate`" + 0.000*`"televis`" + 0.000*`"nowaday`" + 0.000*`"meaning`" + 0.000*`"imag`" + 0.000*`"catch`" + 0.000*`"solv`" + 0.000*`"scarf`" + 0.000*`"bodi`" + 0.000*`"later`" + 0.000*`"entir`" + 0.000*`"display`" + 0.000*`"mental`" + 0.000*`"deserv`" + 0.000*`"straight`" + 0.000*`"minor`" + 0.000*`"tend`" + 0.000*`"weird`" + 0.000*`"hungri`" + 0.000*`"ruin`" + 0.000*`"eye`" + 0.000*`"integr`" + 0.000*`"corner`" + 0.000*`"abandon`" + 0.000*`"expect`" + 0.000*`"woman`" + 0.000*`"disappear`" + 0.000*`"1948`" + 0.000*`"inspir`"')"
$d.Paragraphs.Item(4).Range.Text = "(2, u'0.000*`"intellig`" + 0.000*`"spirit`" + 0.000*`"anj`" + 0.000*`"effort`" + 0.000*`"stranger`" + 0.000*`"unhealthi`" + 0.000*`"onto`" + 0.000*`"brutal`" + 0.000*`"stadium`" + 0.000*`"eljezni`" + 0.000*`"kenan`" + 0.000*`"dread`" + 0.000*`"repair`" + 0.000*`"kahvi`" + 0.000*`"environ`" + 0.000*`"heavili`" + 0.000*`"unlik`" + 0.000*`"master`" + 0.000*`"option`" + 0.000*`"enrol`" + 0.000*`"murder`" + 0.000*`"tool`" + 0.000*`"function`" + 0.000*`"degre`" + 0.000*`"fear`" + 0.000*`"must`" + 0.000*`"keen`" + 0.000*`"diplomaci`" + 0.000*`"jenif`" + 0.000*`"capabl`"')"
$d.Paragraphs.Item(6).Range.Text = "(0, u'0.002*expir + 0.002*lcd + 0.001*far + 0.001*counter + 0.001*agricultur + 0.001*dol + 0.001*enabl + 0.001*steretyp + 0.001*kre + 0.001*architectur + 0.001*symbol + 0.001*progress + 0.001*entri + 0.001*cook + 0.001*ezva + 0.001*technolog + 0.001*dana + 0.001*mcdonald + 0.001*hitler + 0.001*ovo + 0.001*hutu + 0.001*takvu + 0.001*begin + 0.001*nekih + 0.001*die + 0.001*ass + 0.001*last + 0.001*might + 0.001*bought + 0.001*ironi')"
$d.Paragraphs.Item(7).Range.Text = "(1, u'0.003*psychoanalyt + 0.002*felt + 0.001*psihologiju + 0.001*kakav + 0.001*svako + 0.001*krene + 0.001*arthist + 0.001*switch + 0.001*scari + 0.001*distanc + 0.001*auto + 0.001*nephew + 0.001*dobr + 0.001*vrijem + 0.001*10km + 0.001*later + 0.001*bosnish + 0.001*shout + 0.001*beacus + 0.001*preoccupi + 0.001*sell + 0.001*rare + 0.001*ovi + 0.001*hunger + 0.001*organis + 0.001*sunday + 0.001*razlik + 0.001*odem + 0.001*sea + 0.001*nikako')"
$d.Paragraphs.Item(8).Range.Text = "(2, u'0.002*wall + 0.002*began + 0.002*akademiju + 0.002*relax + 0.001*spot + 0.001*dynam + 0.001*2004 + 0.001*divid + 0.001*umrla + 0.001*jer + 0.001*bowel + 0.001*amatur + 0.001*puppi + 0.001*velika + 0.001*tell + 0.001*broj + 0.001*amateur + 0.001*remind + 0.001*psychologist + 0.001*ridicul + 0.001*discov + 0.001*prijatelj + 0.001*conspiraci + 0.001*counselour + 0.001*psihologija + 0.001*vare + 0.001*unlik + 0.001*hot + 0.001*certainli + 0.001*ipad')"
$d.Paragraphs.Item(10).Range.Text = "(0, u'0.002*expir + 0.002*lcd + 0.001*far + 0.001*counter + 0.001*agricultur + 0.001*dol + 0.001*enabl + 0.001*steretyp + 0.001*kre + 0.001*architectur + 0.001*symbol + 0.001*progress + 0.001*entri + 0.001*cook + 0.001*ezva + 0.001*technolog + 0.001*dana + 0.001*mcdonald + 0.001*hitler + 0.001*ovo + 0.001*hutu + 0.001*takvu + 0.001*begin + 0.001*nekih + 0.001*die + 0.001*ass + 0.001*last + 0.001*might + 0.001*bought + 0.001*ironi')"
$d.Paragraphs.Item(11).Range.Text = "(1, u'0.003*psychoanalyt + 0.002*felt + 0.001*psihologiju + 0.001*kakav + 0.001*svako + 0.001*krene + 0.001*arthist + 0.001*switch + 0.001*scari + 0.001*distanc + 0.001*auto + 0.001*nephew + 0.001*dobr + 0.001*vrijem + 0.001*10km + 0.001*later + 0.001*bosnish + 0.001*shout + 0.001*beacus + 0.001*preoccupi + 0.001*sell + 0.001*rare + 0.001*ovi + 0.001*hunger + 0.001*organis + 0.001*sunday + 0.001*razlik + 0.001*odem + 0.001*sea + 0.001*nikako')"
$d.Paragraphs.Item(12).Range.Text = "(2, u'0.002*wall + 0.002*began + 0.002*akademiju + 0.002*relax + 0.001*spot + 0.001*dynam + 0.001*2004 + 0.001*divid + 0.001*umrla + 0.001*jer + 0.001*bowel + 0.001*amatur + 0.001*puppi + 0.001*velika + 0.001*tell + 0.001*broj + 0.001*amateur + 0.001*remind + 0.001*psychologist + 0.001*ridicul + 0.001*discov + 0.001*prijatelj + 0.001*conspiraci + 0.001*counselour + 0.001*psihologija + 0.001*vare + 0.001*unlik + 0.001*hot + 0.001*certainli + 0.001*ipad')"
